$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.253638744354248
$ws.Range("B1").Value = 2.351215124130249
$ws.Range("C1").Value = 1.738563060760498
$ws.Range("D1").Value = 1.620709657669067
$ws.Range("E1").Value = 1.569837093353271
